# Automatic update of files.
# Applies cell-level corrections to rows 5-11 of the "Artfynd" sheet:
# updated Taxonsorteringsordning (B) ids, a re-matching of observation
# records (A/D/E/F/G/H/I/J/P/Q/R/Z/AB) between rows 7-11, and the
# associated coordinate/time fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("B5").Value = 78713

# Row 6
$ws.Range("B6").Value = 89094

# Row 7
$ws.Range("A7").Value = 112212882
$ws.Range("B7").Value = 96735
$ws.Range("J7").Value = "plantor/tuvor"
$ws.Range("P7").Value = "Kälen (Kälen), Jmt"
$ws.Range("Q7").Value = 490109
$ws.Range("R7").Value = 6948768
$ws.Range("Z7").Value = "12:39"
$ws.Range("AB7").Value = "12:39"

# Row 8
$ws.Range("A8").Value = 112212902
$ws.Range("B8").Value = 78713
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6458
$ws.Range("F8").Value = "Lunglav"
$ws.Range("G8").Value = "Lobaria pulmonaria"
$ws.Range("H8").Value = "(L.) Hoffm."
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("Q8").Value = 490134
$ws.Range("R8").Value = 6948772

# Row 9
$ws.Range("A9").Value = 112213279
$ws.Range("B9").Value = 89553
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 1202
$ws.Range("F9").Value = "Ullticka"
$ws.Range("G9").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H9").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I9").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("P9").Value = "Nordvallen (Nordvallen), Jmt"
$ws.Range("Q9").Value = 490080
$ws.Range("R9").Value = 6948907

# Row 10
$ws.Range("A10").Value = 112212105
$ws.Range("B10").Value = 96735
$ws.Range("D10").Value = "VU"
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = "Knärot"
$ws.Range("G10").Value = "Goodyera repens"
$ws.Range("H10").Value = "(L.) R. Br."
# "Antal" (I) is stored as text in this sheet, not a number - use a
# quote-prefix so the numeric-looking "3" is written as text, matching
# how every other populated cell in this column is typed.
$ws.Range("I10").Value = "'3"
$ws.Range("Q10").Value = 490018
$ws.Range("R10").Value = 6948882
$ws.Range("Z10").Value = "11:58"
$ws.Range("AB10").Value = "11:58"

# Row 11
$ws.Range("A11").Value = 112212836
$ws.Range("B11").Value = 96735
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 220787
$ws.Range("F11").Value = "Knärot"
$ws.Range("G11").Value = "Goodyera repens"
$ws.Range("H11").Value = "(L.) R. Br."
$ws.Range("I11").Value = "'25"
$ws.Range("J11").Value = "plantor/tuvor"
$ws.Range("P11").Value = "Stugunäset (Stugunäset), Jmt"
$ws.Range("Q11").Value = 490078
$ws.Range("R11").Value = 6948752
